$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt2"
$ws.Range("C2").Value = "Fzd2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01070233333333333
$ws.Range("H2").Value = 0.032107
$ws.Range("I2").Value = 0.004227647500550067
$ws.Range("J2").Value = 0.004227647500550067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06449866666666666
$ws.Range("N2").Value = 0.193496
$ws.Range("O2").Value = 0.004525829983623641
$ws.Range("P2").Value = 0.004525829983623642
$ws.Range("Q2").Value = 0.0006902862302222222
$ws.Range("R2").Value = 0.006212576072
$ws.Range("S2").Value = 0.00001913361381818104
$ws.Range("T2").Value = 0.00001913361381818104

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt2"
$ws.Range("C3").Value = "Fzd2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01070233333333333
$ws.Range("H3").Value = 0.032107
$ws.Range("I3").Value = 0.004227647500550067
$ws.Range("J3").Value = 0.004227647500550067
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.61985133333333
$ws.Range("N3").Value = 31.859554
$ws.Range("O3").Value = 0.745188142173877
$ws.Range("P3").Value = 0.7451881421738772
$ws.Range("Q3").Value = 0.1136571889197778
$ws.Range("R3").Value = 1.022914700278
$ws.Range("S3").Value = 0.003150392786700939
$ws.Range("T3").Value = 0.00315039278670094

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01070233333333333
$ws.Range("H4").Value = 0.032107
$ws.Range("I4").Value = 0.004227647500550067
$ws.Range("J4").Value = 0.004227647500550067
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.566885000000001
$ws.Range("N4").Value = 10.700655
$ws.Range("O4").Value = 0.2502860278424993
$ws.Range("P4").Value = 0.2502860278424993
$ws.Range("Q4").Value = 0.03817399223166667
$ws.Range("R4").Value = 0.343565930085
$ws.Range("S4").Value = 0.001058121100030946
$ws.Range("T4").Value = 0.001058121100030947

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.520808
$ws.Range("H5").Value = 7.562424
$ws.Range("I5").Value = 0.99577235249945
$ws.Range("J5").Value = 0.99577235249945
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.06449866666666666
$ws.Range("N5").Value = 0.193496
$ws.Range("O5").Value = 0.004525829983623641
$ws.Range("P5").Value = 0.004525829983623642
$ws.Range("Q5").Value = 0.1625887549226667
$ws.Range("R5").Value = 1.463298794304
$ws.Range("S5").Value = 0.00450669636980546
$ws.Range("T5").Value = 0.004506696369805461

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.520808
$ws.Range("H6").Value = 7.562424
$ws.Range("I6").Value = 0.99577235249945
$ws.Range("J6").Value = 0.99577235249945
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.61985133333333
$ws.Range("N6").Value = 31.859554
$ws.Range("O6").Value = 0.745188142173877
$ws.Range("P6").Value = 0.7451881421738772
$ws.Range("Q6").Value = 26.77060619987734
$ws.Range("R6").Value = 240.935455798896
$ws.Range("S6").Value = 0.742037749387176
$ws.Range("T6").Value = 0.7420377493871763

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt2"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.520808
$ws.Range("H7").Value = 7.562424
$ws.Range("I7").Value = 0.99577235249945
$ws.Range("J7").Value = 0.99577235249945
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.566885000000001
$ws.Range("N7").Value = 10.700655
$ws.Range("O7").Value = 0.2502860278424993
$ws.Range("P7").Value = 0.2502860278424993
$ws.Range("Q7").Value = 8.991432243080002
$ws.Range("R7").Value = 80.92289018772001
$ws.Range("S7").Value = 0.2492279067424683
$ws.Range("T7").Value = 0.2492279067424684
